$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. '1.001', '4.170') are preserved exactly instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.761.70'
$ws.Range("E2").Value = '  +1.16%  '
$ws.Range("D3").Value = '1.889.34'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '246.50'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4731'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = '0.2918'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.06522'
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").Value = '22.04'
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").Value = '0.07792'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.890.17'
$ws.Range("E12").Value = '  +1.07%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '96.59'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = '0.7362'
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = '5.241'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("D16").Value = '284.86'
$ws.Range("E16").Value = '  +4.14%  '
$ws.Range("D17").Value = '30.754.48'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("D18").Value = '13.23'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007522'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '2.140.60'
$ws.Range("E21").Value = '  +1.16%  '
$ws.Range("D22").Value = '5.323'
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '6.242'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '9.182'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("D26").Value = '164.46'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").Value = '18.98'
$ws.Range("E27").Value = '  +1.01%  '
$ws.Range("D28").Value = '1.907'
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").Value = '1.338'
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("D30").Value = '0.09709'
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("D31").Value = '1.486'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("D32").Value = '4.293'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '4.170'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").Value = '0.04855'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("D36").Value = '0.6947'
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("D38").Value = '0.01884'
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("D39").Value = '2.802'
$ws.Range("E39").Value = '  +1.91%  '
$ws.Range("D40").Value = '76.18'
$ws.Range("E40").Value = '  +3.72%  '
$ws.Range("D41").Value = '6.295'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '1.990'
$ws.Range("E42").Value = '  +1.48%  '
$ws.Range("D43").Value = '0.4266'
$ws.Range("E43").Value = '  +2.03%  '
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = '0.8341'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").Value = '101.49'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = '9.460'
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").Value = '7.018'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = '35.45'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").Value = '911.54'
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("E51").Value = '  +2.03%  '
